# Add a new action "scroll into view" as a new row in the Actions table (Table1)
# on Sheet1, matching the commit: add action "scroll into view" add lesson for PowerPoint 2010

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row (this also updates the table ref / autofilter / dimension)
$newRow = $lo.ListRows.Add() | Out-Null

# Populate the new row's cells
$ws.Range("A7").Value = "scroll into view"
$ws.Range("B7").Value = "<window_name>|<control_name(list item)>"

# Match the bold style used by the other "Action Name" cells in column A
$ws.Range("A7").Font.Bold = $true

# Move the active selection to the newly added cell, as in the target workbook
$ws.Range("B7").Select() | Out-Null
